$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its original location (just before the
#    "Conclusion" heading run).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Split the big run in the final ("Conclusion") paragraph so that
#    "Artificial Life based" becomes "agent-based", spread across three runs:
#       " is an "  +  "agent-"  +  "based simulation of the ecosystem..."
# ---------------------------------------------------------------------------
$oldText = " is an Artificial Life based simulation of the ecosystem of mobile app development and downloading. The simulation compares different development strategies and shows how these strategies interact to produce popular or unpopular applications. It seems that a Copycat strategy is the most lucrative in terms of downloads, but this strategy must only be used by a minority of developers so that there is a sufficient supply of good apps to be copied. Future research is planned to investigate the effects of publicity on app downloads, and explore how users might be able to better locate desirable apps and provide feedback of their preferences back to developers."

$r = $d.Content
$found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $r.Delete()

    $part1 = " is an "
    $part2 = "agent-"
    $part3 = "based simulation of the ecosystem of mobile app development and downloading. The simulation compares different development strategies and shows how these strategies interact to produce popular or unpopular applications. It seems that a Copycat strategy is the most lucrative in terms of downloads, but this strategy must only be used by a minority of developers so that there is a sufficient supply of good apps to be copied. Future research is planned to investigate the effects of publicity on app downloads, and explore how users might be able to better locate desirable apps and provide feedback of their preferences back to developers."

    $r.InsertAfter($part1)
    $r.Collapse(0)

    $r.InsertAfter($part2)
    $r.Collapse(0)

    $r.InsertAfter($part3)
    $r.Collapse(0)
}

# ---------------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark, collapsed, at the very end of the document
#    (after the last run of the final paragraph).
#
#    The engine this runs on fails to create a bookmark collapsed exactly at
#    the true end of the document content, so we temporarily pad the document
#    with a few extra characters, add the bookmark at the (now not-final)
#    position, and then remove the padding again.
# ---------------------------------------------------------------------------
$rawEnd = $d.Content.End
$realEnd = $rawEnd - 1

$padRange = $d.Range($rawEnd, $rawEnd)
$padRange.InsertAfter("PAD")

$bmRange = $d.Range($realEnd, $realEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padDelRange = $d.Range($realEnd, $realEnd + 3)
$padDelRange.Delete()
